$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8,1).Value = 10
$ws.Cells.Item(8,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8,3).Value = "La Araucanía"
$ws.Cells.Item(8,4).Value = 44503
$ws.Cells.Item(8,5).Value = 9
$ws.Cells.Item(8,6).Value = 100112022
$ws.Cells.Item(8,7).Value = "Arveja Verde"
$ws.Cells.Item(8,8).Value = "Sin especificar"
$ws.Cells.Item(8,9).Value = "Primera"
$ws.Cells.Item(8,10).Value = 75
$ws.Cells.Item(8,11).Value = 15000
$ws.Cells.Item(8,12).Value = 15000
$ws.Cells.Item(8,13).Value = 15000
$ws.Cells.Item(8,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(8,15).Value = "Provincia de Limarí"
$ws.Cells.Item(8,16).Value = 600
$ws.Cells.Item(8,17).Value = 25
$ws.Cells.Item(8,18).Value = "Hortaliza"
